# Applies the red_phash.xlsx edit: updates several "Время обработки" (column C)
# values on existing rows, and appends 20 new rows (red_fred_1.jpg .. red_fred_20.jpg)
# after the existing data, extending the sheet from A1:D64 to A1:D84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update column C ("Время обработки") values on existing rows 2-61.
# ---------------------------------------------------------------------------
$cUpdates = @(
    @{Row=2;  Value=0},
    @{Row=3;  Value=0},
    @{Row=4;  Value=0.015614},
    @{Row=15; Value=0.01562},
    @{Row=16; Value=0},
    @{Row=29; Value=0.015622},
    @{Row=30; Value=0},
    @{Row=42; Value=0.015627},
    @{Row=51; Value=0.01562},
    @{Row=56; Value=0.015639},
    @{Row=61; Value=0.015615}
)

foreach ($u in $cUpdates) {
    $ws.Cells.Item($u.Row, 3).Value2 = $u.Value
}

# ---------------------------------------------------------------------------
# 2) Append 20 new rows (65-84) with the red_fred_* data.
# ---------------------------------------------------------------------------
$hashB = "1110101000000000110000001000000010000110000000001000000000000000"
$hashB19 = "1110101000000000110000000000000010000110000000001000000000000000"

$newRows = @(
    @{Row=65; A="red_fred_1.jpg";  B=$hashB;   C=0;        D=0},
    @{Row=66; A="red_fred_2.jpg";  B=$hashB;   C=0.015627; D=0},
    @{Row=67; A="red_fred_3.jpg";  B=$hashB;   C=0;        D=0},
    @{Row=68; A="red_fred_4.jpg";  B=$hashB;   C=0;        D=0},
    @{Row=69; A="red_fred_5.jpg";  B=$hashB;   C=0;        D=0},
    @{Row=70; A="red_fred_6.jpg";  B=$hashB;   C=0;        D=0},
    @{Row=71; A="red_fred_7.jpg";  B=$hashB;   C=0.015622; D=0},
    @{Row=72; A="red_fred_8.jpg";  B=$hashB;   C=0;        D=0},
    @{Row=73; A="red_fred_9.jpg";  B=$hashB;   C=0;        D=0},
    @{Row=74; A="red_fred_10.jpg"; B=$hashB;   C=0;        D=0},
    @{Row=75; A="red_fred_11.jpg"; B=$hashB;   C=0;        D=0},
    @{Row=76; A="red_fred_12.jpg"; B=$hashB;   C=0.015628; D=0},
    @{Row=77; A="red_fred_13.jpg"; B=$hashB;   C=0;        D=0},
    @{Row=78; A="red_fred_14.jpg"; B=$hashB;   C=0;        D=0},
    @{Row=79; A="red_fred_15.jpg"; B=$hashB;   C=0;        D=0},
    @{Row=80; A="red_fred_16.jpg"; B=$hashB;   C=0;        D=0},
    @{Row=81; A="red_fred_17.jpg"; B=$hashB;   C=0.015633; D=0},
    @{Row=82; A="red_fred_18.jpg"; B=$hashB;   C=0;        D=0},
    @{Row=83; A="red_fred_19.jpg"; B=$hashB19; C=0;        D=1},
    @{Row=84; A="red_fred_20.jpg"; B=$hashB;   C=0;        D=0}
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value2 = $r.A

    # Column B holds a long string of digits ("0"/"1" bits). Excel's COM
    # layer auto-coerces a pure-digit string into a number (and then into
    # scientific notation) unless the cell is pre-formatted as Text, so we
    # force a text format before assigning the value.
    $cellB = $ws.Cells.Item($r.Row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value2 = $r.B

    $ws.Cells.Item($r.Row, 3).Value2 = $r.C
    $ws.Cells.Item($r.Row, 4).Value2 = $r.D
}
